# Apply the "BOs report" refresh: new occurrence row data, shrink the
# autofilter/_FilterDatabase range down to just the header row, even out
# the widths of columns Q:R, and drop the trailing block of empty rows
# that used to pad the sheet out to row 86.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Refresh the single data row (row 2) with the new occurrence record ---
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = 177745
$ws.Range("C2").Value = "laion de Lara silva"
$ws.Range("D2").Value = "UDT-1E38"
$ws.Range("F2").Value = "2026-01-16 04:04:57.000000 UTC"
$ws.Range("G2").Value = "-23.975288, -48.864197"
$ws.Range("J2").Value = "R. Higino Marques, 87, Itapeva - SP, 18407120"
$ws.Range("K2").Value = "2026-01-15 04:41:18.020000 UTC"
$ws.Range("L2").Value = "-23.975288, -48.864197"
$ws.Range("M2").Value = 4750020
$ws.Range("N2").Value = 4429466
$ws.Range("O2").Value = "-"
$ws.Range("P2").Value = "92EC10BHSSM055161"
$ws.Range("Q2").Value = "Mottu Itapeva"

# --- 2. Even out column R's width so it matches column Q ---
$ws.Columns.Item(18).ColumnWidth = $ws.Columns.Item(17).ColumnWidth()

# --- 3. Drop the 22 trailing blank rows (65-86) that used to pad the sheet ---
$ws.Range("A65:A86").EntireRow.Delete()

# --- 4. Shrink the autofilter down to the header row only (A1:T1) ---
# (Range.AutoFilter with no args toggles filtering for the given range, so
#  call it twice: once to turn the existing filter off, once to turn it
#  back on against the new A1:T1 range.)
[void]$ws.Range("A1:T1").AutoFilter()
[void]$ws.Range("A1:T1").AutoFilter()

# --- 5. Keep the hidden _FilterDatabase defined name in sync with the autofilter ---
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='" + $ws.Name + "'!`$A`$1:`$T`$1"
    }
}
